$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd name "Prretika Shetty" -> "Preetika Shetty"
# (this also moves the string to the end of the shared-strings table
# because the cell referencing it, B3, is rewritten below)
$ws.Range("B3").Value = "Preetika Shetty"

# Update column A (id numbers) and column C (scores) for all 10 rows
$ws.Range("A1").Value = 52501
$ws.Range("A2").Value = 52502
$ws.Range("A3").Value = 52503
$ws.Range("A4").Value = 52504
$ws.Range("A5").Value = 52505
$ws.Range("A6").Value = 52506
$ws.Range("A7").Value = 52507
$ws.Range("A8").Value = 52508
$ws.Range("A9").Value = 52509
$ws.Range("A10").Value = 52510

$ws.Range("C1").Value = 80
$ws.Range("C2").Value = 80
$ws.Range("C3").Value = 97
$ws.Range("C4").Value = 96
$ws.Range("C5").Value = 93
$ws.Range("C6").Value = 95
$ws.Range("C7").Value = 92
$ws.Range("C8").Value = 85
$ws.Range("C9").Value = 91
$ws.Range("C10").Value = 90

# Update the selection to match the saved view state (active cell B3)
$ws.Range("B3").Select()
